# "process macro v4.3 added"
#
# The "Video Links" document originally listed four video sections
# (Three types of missing data / Linear Regression / Multiple Linear
# Regression / Linear Regression Summary in R), each followed by its
# YouTube hyperlink. This edit:
#   1. Retitles the first heading to "Installing PROCESS Macro v4.3".
#   2. Removes that heading's old hyperlink paragraph.
#   3. Removes every other video section (heading + hyperlink + spacer)
#      entirely, leaving just the new heading and the trailing blank
#      paragraphs.
#   4. Best-effort: marks the "Default Paragraph Font" style as hidden
#      in the style gallery (maps to <w:semiHidden/> in styles.xml).

$d = $word.ActiveDocument

# 1. Retitle the first heading in place (keeps its bold run formatting).
$d.Content.Find.Execute(
    "Three types of missing data " + [char]0x2013 + " 3 minutes",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Installing PROCESS Macro v4.3", 2) | Out-Null

# 2. Remove the hyperlink paragraph that used to follow that heading
#    (paragraph 2 after the rename: the "UzsWr9X98J8" link), including
#    its paragraph mark, leaving the blank paragraph after it intact.
$r = $d.Range($d.Paragraphs(2).Range.Start, $d.Paragraphs(2).Range.End)
$r.Delete()

# 3. Remove everything from "Linear Regression clearly explained..."
#    through the final "7WPfuHLCn_k" hyperlink paragraph (now
#    paragraphs 3-10), including their paragraph marks. This drops the
#    Linear Regression / Multiple Linear Regression / Linear Regression
#    Summary in R sections and their hyperlinks and spacer paragraphs,
#    while leaving the two trailing blank paragraphs untouched.
$r2 = $d.Range($d.Paragraphs(3).Range.Start, $d.Paragraphs(10).Range.End)
$r2.Delete()

# 4. Best-effort style tweak: hide "Default Paragraph Font" from the
#    style gallery (<w:semiHidden/>). Not all runtimes expose a writable
#    setter for this; ignore failures so the content edits above stand
#    regardless.
try {
    $style = $d.Styles("Default Paragraph Font")
    $style.Hidden = $true
} catch {
}
